# Sort the data by project_year ascending, then by country_impact_name ascending
# (per the commit message: "added sorting by ascending order per year")
#
# Note: Excel's Range.Sort always pushes blank cells to the end of an ascending
# sort, regardless of sort order. The target layout instead expects blank
# country names to sort first (as the empty string sorts before any letter).
# To achieve that we temporarily replace blanks with a sentinel value that
# collates before "A", sort, then clear the sentinel back to blank.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Determine the extent of the data (header in row 1, data starts row 2)
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row   # xlUp = -4162
$lastCol = 2  # columns A (project_year) and B (country_impact_name)

$dataRange = $ws.Range($ws.Cells.Item(1,1), $ws.Cells.Item($lastRow,$lastCol))

$keyYear = $ws.Range($ws.Cells.Item(1,1), $ws.Cells.Item($lastRow,1))
$keyCountry = $ws.Range($ws.Cells.Item(1,2), $ws.Cells.Item($lastRow,2))

# Sentinel that sorts before any real country name
$sentinel = "!!!BLANK!!!"

# Remember which rows in column B are blank, then fill them with the sentinel
$blankRows = New-Object System.Collections.ArrayList
for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 2)
    if ([string]::IsNullOrEmpty($cell.Value())) {
        $blankRows.Add($r) | Out-Null
        $cell.Value = $sentinel
    }
}

$ws.Sort.SortFields.Clear()
$ws.Sort.SortFields.Add($keyYear, 0, 1, 0, 0) | Out-Null
$ws.Sort.SortFields.Add($keyCountry, 0, 1, 0, 0) | Out-Null

$ws.Sort.SetRange($dataRange)
$ws.Sort.Header = 1        # xlYes : first row is header, excluded from sort
$ws.Sort.MatchCase = $false
$ws.Sort.Orientation = 1   # xlTopToBottom
$ws.Sort.SortMethod = 1    # xlPinYin
$ws.Sort.Apply()

# Restore the sentinel cells back to blank
for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 2)
    if ($cell.Value() -eq $sentinel) {
        $cell.Value = ""
    }
}
